# Refresh this NATMI ligand-receptor (Fn1 -> Itgb3) sheet with the new TPM
# expression values. The raw per-cluster ligand (G/H) and receptor
# (K/L/M/N) numbers come from the updated TPM table; every other numeric
# column on the sheet (I/J/O/P specificities and Q/R/S/T edge weights) is
# derived from those five-cluster-wide raw numbers, so we recompute them
# here instead of poking 288 cells by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand expression (per sending cluster), from the refreshed TPM data.
$ligandAvg = @{
    "ECs"               = 66.47695399999999
    "FAPs"              = 1361.379069
    "Inflammatory-Mac"  = 44.831112
    "MuSCs"             = 52.83062100000001
    "Resolving-Mac"     = 16.16161433333333
}
$ligandTotal = @{
    "ECs"               = 199.430862
    "FAPs"              = 4084.137207
    "Inflammatory-Mac"  = 134.493336
    "MuSCs"             = 158.491863
    "Resolving-Mac"     = 48.484843
}

# New receptor expression (per target cluster), from the refreshed TPM data.
$receptorCells = @{
    "ECs"               = 3
    "FAPs"              = 3
    "Inflammatory-Mac"  = 3
    "MuSCs"             = 3
    "Resolving-Mac"     = 3
}
$receptorRate = @{
    "ECs"               = 1
    "FAPs"              = 1
    "Inflammatory-Mac"  = 1
    "MuSCs"             = 1
    "Resolving-Mac"     = 1
}
$receptorAvg = @{
    "ECs"               = 6.111751666666666
    "FAPs"              = 2.754304
    "Inflammatory-Mac"  = 0.568439
    "MuSCs"             = 0.3689163333333333
    "Resolving-Mac"     = 0.2794673333333333
}
$receptorTotal = @{
    "ECs"               = 18.335255
    "FAPs"              = 8.262912
    "Inflammatory-Mac"  = 1.705317
    "MuSCs"             = 1.106749
    "Resolving-Mac"     = 0.838402
}

# Specificity = a cluster's value over the sum across all clusters.
$sumLigandAvg = 0
$sumLigandTotal = 0
$sumReceptorAvg = 0
$sumReceptorTotal = 0
foreach ($cluster in $ligandAvg.Keys) {
    $sumLigandAvg += $ligandAvg[$cluster]
    $sumLigandTotal += $ligandTotal[$cluster]
    $sumReceptorAvg += $receptorAvg[$cluster]
    $sumReceptorTotal += $receptorTotal[$cluster]
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 26 }

for ($row = 2; $row -le $lastRow; $row++) {
    $sendCluster = $ws.Cells.Item($row, 1).Value()
    $targetCluster = $ws.Cells.Item($row, 4).Value()
    if (-not $sendCluster -or -not $ligandAvg.ContainsKey($sendCluster)) { continue }
    if (-not $targetCluster -or -not $receptorAvg.ContainsKey($targetCluster)) { continue }

    $g = $ligandAvg[$sendCluster]
    $h = $ligandTotal[$sendCluster]
    $i = $g / $sumLigandAvg
    $j = $h / $sumLigandTotal

    $k = $receptorCells[$targetCluster]
    $l = $receptorRate[$targetCluster]
    $m = $receptorAvg[$targetCluster]
    $n = $receptorTotal[$targetCluster]
    $o = $m / $sumReceptorAvg
    $p = $n / $sumReceptorTotal

    $q = $g * $m
    $r = $h * $n
    $s = $i * $o
    $t = $j * $p

    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $j
    $ws.Cells.Item($row, 11).Value = $k
    $ws.Cells.Item($row, 12).Value = $l
    $ws.Cells.Item($row, 13).Value = $m
    $ws.Cells.Item($row, 14).Value = $n
    $ws.Cells.Item($row, 15).Value = $o
    $ws.Cells.Item($row, 16).Value = $p
    $ws.Cells.Item($row, 17).Value = $q
    $ws.Cells.Item($row, 18).Value = $r
    $ws.Cells.Item($row, 19).Value = $s
    $ws.Cells.Item($row, 20).Value = $t
}
